$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row before row 8 (existing rows 8.. shift down by one) ---
$ws.Rows.Item(8).Insert()

# Copy the formatting (styles: s="3"/s="4"/s="3"/s="3") from row 7, which
# already carries the same B/C/D/E style pattern the new row needs.
$ws.Range("B7:E7").Copy()
$ws.Range("B8:E8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Fill in the new row's content: "json" type + its restriction pattern ---
$ws.Cells.Item(8, 2).Value = "json"
$ws.Cells.Item(8, 3).Value = "`"type`": [ `"string`" ]," + [char]10 + "`"format`": `"string`"," + [char]10 + "`"pattern`": `"^(([{](((\`"[a-zA-Z0-9]{1,}\`":((null)|(([0-9]{1,})((.)([0-9]{1,})){0,1})|(\`".{1,}\`")|([[].*[]])|([{]\`"[a-zA-Z0-9]{1,}\`":(.*)[}]))),){0,}(\`"[a-zA-Z0-9]{1,}\`":((null)|([0-9]{1,})|(\`"[a-zA-Z0-9]{1,}\`")|([[].*[]])|([{]\`"[a-zA-Z0-9]{1,}\`":(.*)[}]))){0,})[}])|([[]((((null)|(([0-9]{1,})((.)([0-9]{1,})){0,1})|(\`".{1,}\`")|([[].*[]])|([{]\`"[a-zA-Z0-9]{1,}\`":(.*)[}])),){0,}((null)|(([0-9]{1,})((.)([0-9]{1,})){0,1})|(\`".{1,}\`")|([[].*[]])|([{]\`"[a-zA-Z0-9]{1,}\`":(.*)[}])))[]]))`$`""

# Row 8 needs the taller, 5-line-tall row height (its new pattern text wraps
# further than the row it was cloned from).
$ws.Rows.Item(8).RowHeight = 63.75

# --- Restore the frozen-pane split (Insert/PasteSpecial can disturb it) and
#     move the view's active cell down to account for the inserted row ---
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Application.Goto($ws.Range("C3"))
$win.FreezePanes = $true
$ws.Application.Goto($ws.Range("C9"))
